$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix reviewer / product name in the second (user documentation) table ---
# Order matters for shared-string allocation: reviewer name first, then the
# corrected product title.
$ws.Range("B11").Value = "Krizsák Kornél"
$ws.Range("B10").Value = "IKT_Webshop Felhasználói dokumentáció"

# --- Bug #1 (row 14): popup window closes on outside click too ---
$ws.Range("C14").Value = "3. oldal/5-7"
$ws.Range("B14").Value = "A szöveg eredetileg a felugró ablak bezárására csak az erre megadott gombbal kínál lehetőséget, valójában mellékattintva is bezárul az ablak."
$ws.Range("D14").Value = "jelentékeny"

# --- Bug #2 (row 15): ads are not actually random ---
$ws.Range("C15").Value = "3. oldal/14"
$ws.Range("B15").Value = "A szöveg véletlenszerű reklámra hivatkozik, miközben a reklámok előre megadottak"
$ws.Range("D15").Value = "elhanyagolható"

# --- Bug #3 (row 16): comment info is limited, not detailed ---
$ws.Range("C16").Value = "3. oldal/20"
$ws.Range("B16").Value = "A szöveg a megjegyzés részletes információira hivatkozik, valójában csak korlátozott információk elérhetőek"
$ws.Range("D16").Value = "jelentékeny"

# Long descriptions need to wrap, and their rows grow taller to fit.
$ws.Range("B14").WrapText = $true
$ws.Range("B16").WrapText = $true
$ws.Rows.Item(14).RowHeight = 30
$ws.Rows.Item(16).RowHeight = 30

# Final cursor position left by the author before saving.
$ws.Range("G20").Select()
